$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.158.72"
$ws.Range("E2").Value = '  +7.94%  '

$ws.Range("D3").Value = "'1.585.45"
$ws.Range("E3").Value = '  +7.70%  '

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = "'0.9951"
$ws.Range("E5").Value = '  +3.99%  '

$ws.Range("D6").Value = "'298.25"
$ws.Range("E6").Value = '  +7.43%  '

$ws.Range("D7").Value = "'0.3612"
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = "'0.3337"
$ws.Range("E8").Value = '  +8.90%  '

$ws.Range("D9").Value = "'41.26"
$ws.Range("E9").Value = '  +4.22%  '

$ws.Range("E10").Value = '  +4.18%  '

$ws.Range("D11").Value = "'0.06927"
$ws.Range("E11").Value = '  +4.41%  '

$ws.Range("D12").Value = "'1.007"
$ws.Range("E12").Value = '  +0.51%  '

$ws.Range("D13").Value = "'19.24"
$ws.Range("E13").Value = '  +6.39%  '

$ws.Range("D14").Value = "'5.794"
$ws.Range("E14").Value = '  +4.88%  '

$ws.Range("D15").Value = "'6.512"
$ws.Range("E15").Value = '  +5.57%  '

$ws.Range("D16").Value = "'0.9965"
$ws.Range("E16").Value = '  +3.95%  '

$ws.Range("D17").Value = "'1.583.02"
$ws.Range("E17").Value = '  +7.29%  '

$ws.Range("D18").Value = "'0.00001056"

$ws.Range("D19").Value = "'0.06565"
$ws.Range("E19").Value = '  +10.94%  '

$ws.Range("D20").Value = "'75.86"
$ws.Range("E20").Value = '  +9.86%  '

$ws.Range("D21").Value = "'5.895"
$ws.Range("E21").Value = '  +7.33%  '

$ws.Range("E22").Value = '  +8.63%  '

$ws.Range("D23").Value = "'11.58"
$ws.Range("E23").Value = '  +3.77%  '

$ws.Range("D24").Value = "'22.206.96"
$ws.Range("E24").Value = '  +8.07%  '

$ws.Range("D25").Value = "'2.375"
$ws.Range("E25").Value = '  +5.15%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = "'2.467"
$ws.Range("E26").Value = '  +16.38%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'147.81"
$ws.Range("E27").Value = '  +2.95%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'19.08"
$ws.Range("E28").Value = '  +11.19%  '

$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = "'1.753.76"
$ws.Range("E29").Value = '  +7.12%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = "'120.72"
$ws.Range("E30").Value = '  +6.24%  '

$ws.Range("B31").Value = 'HuobiToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D31").Value = "'3.943"
$ws.Range("E31").Value = '  +0.67%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.833"
$ws.Range("E32").Value = '  +17.54%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = "'0.9100"
$ws.Range("E33").Value = '  +12.72%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = "'0.08110"
$ws.Range("E34").Value = '  +1.48%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = "'1.630"
$ws.Range("E35").Value = '  +7.58%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = "'11.61"
$ws.Range("E36").Value = '  +12.09%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = "'5.082"
$ws.Range("E37").Value = '  +7.63%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'1.223"
$ws.Range("E38").Value = '  +0.85%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'8.292"
$ws.Range("E39").Value = '  +11.65%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.05948"
$ws.Range("E40").Value = '  +3.49%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = "'0.02169"
$ws.Range("E41").Value = '  +5.62%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = "'0.9960"
$ws.Range("E42").Value = '  +3.90%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = "'0.1973"
$ws.Range("E43").Value = '  +5.14%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = "'0.5754"
$ws.Range("E44").Value = '  +9.20%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = "'3.761"
$ws.Range("E45").Value = '  +6.72%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'12.83"
$ws.Range("E46").Value = '  +5.03%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = "'124.65"
$ws.Range("E47").Value = '  +5.69%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = "'0.5536"
$ws.Range("E48").Value = '  +6.54%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'1.927"
$ws.Range("E49").Value = '  +6.35%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.06720"
$ws.Range("E50").Value = '  +3.94%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'72.01"
$ws.Range("E51").Value = '  +7.28%  '
